$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

# --- Add new row 8: TestCase7 ---
$ws.Cells.Item(8, 1).Value = "TestCase7"
$ws.Cells.Item(8, 2).Value = "D:\SIN_ADE.pdf"
$ws.Cells.Item(8, 3).Value = "D:\SIN_UW.pdf"
$ws.Cells.Item(8, 4).Value = "No"
$ws.Range("A8:D8").WrapText = $true

# --- Add new row 9: TestCase8 ---
$ws.Cells.Item(9, 1).Value = "TestCase8"
$ws.Cells.Item(9, 2).Value = "D:\Test.pdf"
$ws.Cells.Item(9, 3).Value = "D:\Test.pdf"
$ws.Cells.Item(9, 4).Value = "No"
$ws.Cells.Item(9, 1).WrapText = $true
$ws.Cells.Item(9, 4).WrapText = $true

# --- Add new row 10: TestCase9 ---
$ws.Cells.Item(10, 1).Value = "TestCase9"
$ws.Cells.Item(10, 2).Value = "D:\Test.pdf"
$ws.Cells.Item(10, 3).Value = "D:\Test.pdf"
$ws.Cells.Item(10, 4).Value = "No"
$ws.Cells.Item(10, 1).WrapText = $true
$ws.Cells.Item(10, 4).WrapText = $true

# --- Add new row 11: TestCase10 ---
$ws.Cells.Item(11, 1).Value = "TestCase10"
$ws.Cells.Item(11, 2).Value = "D:\Test.pdf"
$ws.Cells.Item(11, 3).Value = "D:\Test.pdf"
$ws.Cells.Item(11, 4).Value = "No"
$ws.Cells.Item(11, 1).WrapText = $true
$ws.Cells.Item(11, 4).WrapText = $true

# --- Add new row 12: TestCase11 ---
$ws.Cells.Item(12, 1).Value = "TestCase11"
$ws.Cells.Item(12, 2).Value = "D:\Test.pdf"
$ws.Cells.Item(12, 3).Value = "D:\Test.pdf"
$ws.Cells.Item(12, 4).Value = "No"
$ws.Cells.Item(12, 1).WrapText = $true
$ws.Cells.Item(12, 4).WrapText = $true

# --- Update existing row 2 (TestCase1) input files ---
$ws.Cells.Item(2, 2).Value = "D:\ABRCIR-20180326.pdf"
$ws.Cells.Item(2, 3).Value = "D:\Correspondence March 22 (003).pdf"

# --- Flip ExecutionMode from Yes to No for TestCase2..TestCase6 (rows 3-7) ---
$ws.Cells.Item(3, 4).Value = "No"
$ws.Cells.Item(4, 4).Value = "No"
$ws.Cells.Item(5, 4).Value = "No"
$ws.Cells.Item(6, 4).Value = "No"
$ws.Cells.Item(7, 4).Value = "No"

# --- Update selection to match the authored state ---
[void]$ws.Activate()
[void]$ws.Range("C3").Select()

# --- Best-effort: disable concurrent/multi-threaded calculation ---
$excel.MultiThreadedCalculation.Enabled = $false
